$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add consump to stock" -- fill in the previously-empty quarterly "سرمایه"
# (capital/stock) values on row 26 so they match the prior quarters.
$ws.Range("I26").Value = 2000000
$ws.Range("J26").Value = 2000000
$ws.Range("M26").Value = 2000000

# Restore the view state captured when the change was made: right-to-left
# sheet layout, scrolled so column F is at the left edge, with J31 selected.
$excel.ActiveWindow.DisplayRightToLeft = $true
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("J31").Select() | Out-Null
